$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.58598166666667
$ws.Range("H2").Value = 82.75794500000001
$ws.Range("I2").Value = 0.2704460545904799
$ws.Range("J2").Value = 0.2704460545904799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.839587
$ws.Range("N2").Value = 68.518761
$ws.Range("O2").Value = 0.2024156068965367
$ws.Range("P2").Value = 0.2024156068965367
$ws.Range("Q2").Value = 630.0524282562384
$ws.Range("R2").Value = 5670.471854306145
$ws.Range("S2").Value = 0.05474250227270588
$ws.Range("T2").Value = 0.05474250227270588

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.58598166666667
$ws.Range("H3").Value = 82.75794500000001
$ws.Range("I3").Value = 0.2704460545904799
$ws.Range("J3").Value = 0.2704460545904799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 31.29092
$ws.Range("N3").Value = 93.87276
$ws.Range("O3").Value = 0.2773154594323872
$ws.Range("P3").Value = 0.2773154594323872
$ws.Range("Q3").Value = 863.1907454531334
$ws.Range("R3").Value = 7768.7167090782
$ws.Range("S3").Value = 0.0749988718804354
$ws.Range("T3").Value = 0.0749988718804354

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.58598166666667
$ws.Range("H4").Value = 82.75794500000001
$ws.Range("I4").Value = 0.2704460545904799
$ws.Range("J4").Value = 0.2704460545904799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 35.824351
$ws.Range("N4").Value = 107.473053
$ws.Range("O4").Value = 0.3174929454433458
$ws.Range("P4").Value = 0.3174929454433459
$ws.Range("Q4").Value = 988.2498899062317
$ws.Range("R4").Value = 8894.249009156087
$ws.Range("S4").Value = 0.08586471445546337
$ws.Range("T4").Value = 0.08586471445546337

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.58598166666667
$ws.Range("H5").Value = 82.75794500000001
$ws.Range("I5").Value = 0.2704460545904799
$ws.Range("J5").Value = 0.2704460545904799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.88025066666667
$ws.Range("N5").Value = 68.64075199999999
$ws.Range("O5").Value = 0.2027759882277303
$ws.Range("P5").Value = 0.2027759882277303
$ws.Range("Q5").Value = 631.1741754194045
$ws.Range("R5").Value = 5680.56757877464
$ws.Range("S5").Value = 0.05483996598187526
$ws.Range("T5").Value = 0.05483996598187525

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 26.23504533333333
$ws.Range("H6").Value = 78.705136
$ws.Range("I6").Value = 0.2572018131577233
$ws.Range("J6").Value = 0.2572018131577233
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 22.839587
$ws.Range("N6").Value = 68.518761
$ws.Range("O6").Value = 0.2024156068965367
$ws.Range("P6").Value = 0.2024156068965367
$ws.Range("Q6").Value = 599.1976003396106
$ws.Range("R6").Value = 5392.778403056495
$ws.Range("S6").Value = 0.0520616611052102
$ws.Range("T6").Value = 0.0520616611052102

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 26.23504533333333
$ws.Range("H7").Value = 78.705136
$ws.Range("I7").Value = 0.2572018131577233
$ws.Range("J7").Value = 0.2572018131577233
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.29092
$ws.Range("N7").Value = 93.87276
$ws.Range("O7").Value = 0.2773154594323872
$ws.Range("P7").Value = 0.2773154594323872
$ws.Range("Q7").Value = 820.9187047217066
$ws.Range("R7").Value = 7388.268342495359
$ws.Range("S7").Value = 0.07132603898267705
$ws.Range("T7").Value = 0.07132603898267705

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.23504533333333
$ws.Range("H8").Value = 78.705136
$ws.Range("I8").Value = 0.2572018131577233
$ws.Range("J8").Value = 0.2572018131577233
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 35.824351
$ws.Range("N8").Value = 107.473053
$ws.Range("O8").Value = 0.3174929454433458
$ws.Range("P8").Value = 0.3174929454433459
$ws.Range("Q8").Value = 939.8534725222453
$ws.Range("R8").Value = 8458.681252700208
$ws.Range("S8").Value = 0.08165976123281468
$ws.Range("T8").Value = 0.08165976123281468

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.23504533333333
$ws.Range("H9").Value = 78.705136
$ws.Range("I9").Value = 0.2572018131577233
$ws.Range("J9").Value = 0.2572018131577233
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.88025066666667
$ws.Range("N9").Value = 68.64075199999999
$ws.Range("O9").Value = 0.2027759882277303
$ws.Range("P9").Value = 0.2027759882277303
$ws.Range("Q9").Value = 600.2644134780302
$ws.Range("R9").Value = 5402.379721302271
$ws.Range("S9").Value = 0.05215435183702138
$ws.Range("T9").Value = 0.05215435183702138

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.31506333333334
$ws.Range("H10").Value = 87.94519000000001
$ws.Range("I10").Value = 0.2873975381543141
$ws.Range("J10").Value = 0.2873975381543141
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 22.839587
$ws.Range("N10").Value = 68.518761
$ws.Range("O10").Value = 0.2024156068965367
$ws.Range("P10").Value = 0.2024156068965367
$ws.Range("Q10").Value = 669.5439394121768
$ws.Range("R10").Value = 6025.895454709591
$ws.Range("S10").Value = 0.05817374710607605
$ws.Range("T10").Value = 0.05817374710607606

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 29.31506333333334
$ws.Range("H11").Value = 87.94519000000001
$ws.Range("I11").Value = 0.2873975381543141
$ws.Range("J11").Value = 0.2873975381543141
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 31.29092
$ws.Range("N11").Value = 93.87276
$ws.Range("O11").Value = 0.2773154594323872
$ws.Range("P11").Value = 0.2773154594323872
$ws.Range("Q11").Value = 917.2953015582668
$ws.Range("R11").Value = 8255.657714024401
$ws.Range("S11").Value = 0.07969978033300064
$ws.Range("T11").Value = 0.07969978033300065

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 29.31506333333334
$ws.Range("H12").Value = 87.94519000000001
$ws.Range("I12").Value = 0.2873975381543141
$ws.Range("J12").Value = 0.2873975381543141
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 35.824351
$ws.Range("N12").Value = 107.473053
$ws.Range("O12").Value = 0.3174929454433458
$ws.Range("P12").Value = 0.3174929454433459
$ws.Range("Q12").Value = 1050.193118440563
$ws.Range("R12").Value = 9451.738065965072
$ws.Range("S12").Value = 0.09124669090177956
$ws.Range("T12").Value = 0.09124669090177957

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 29.31506333333334
$ws.Range("H13").Value = 87.94519000000001
$ws.Range("I13").Value = 0.2873975381543141
$ws.Range("J13").Value = 0.2873975381543141
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 22.88025066666667
$ws.Range("N13").Value = 68.64075199999999
$ws.Range("O13").Value = 0.2027759882277303
$ws.Range("P13").Value = 0.2027759882277303
$ws.Range("Q13").Value = 670.7359973758756
$ws.Range("R13").Value = 6036.62397638288
$ws.Range("S13").Value = 0.05827731981345786
$ws.Range("T13").Value = 0.05827731981345786

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 18.86569966666666
$ws.Range("H14").Value = 56.59709899999999
$ws.Range("I14").Value = 0.1849545940974826
$ws.Range("J14").Value = 0.1849545940974826
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 22.839587
$ws.Range("N14").Value = 68.518761
$ws.Range("O14").Value = 0.2024156068965367
$ws.Range("P14").Value = 0.2024156068965367
$ws.Range("Q14").Value = 430.8847888527042
$ws.Range("R14").Value = 3877.963099674339
$ws.Range("S14").Value = 0.03743769641254455
$ws.Range("T14").Value = 0.03743769641254455

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 18.86569966666666
$ws.Range("H15").Value = 56.59709899999999
$ws.Range("I15").Value = 0.1849545940974826
$ws.Range("J15").Value = 0.1849545940974826
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 31.29092
$ws.Range("N15").Value = 93.87276
$ws.Range("O15").Value = 0.2773154594323872
$ws.Range("P15").Value = 0.2773154594323872
$ws.Range("Q15").Value = 590.3250990136933
$ws.Range("R15").Value = 5312.925891123239
$ws.Range("S15").Value = 0.05129076823627408
$ws.Range("T15").Value = 0.05129076823627409

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 18.86569966666666
$ws.Range("H16").Value = 56.59709899999999
$ws.Range("I16").Value = 0.1849545940974826
$ws.Range("J16").Value = 0.1849545940974826
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 35.824351
$ws.Range("N16").Value = 107.473053
$ws.Range("O16").Value = 0.3174929454433458
$ws.Range("P16").Value = 0.3174929454433459
$ws.Range("Q16").Value = 675.8514467192496
$ws.Range("R16").Value = 6082.663020473246
$ws.Range("S16").Value = 0.05872177885328823
$ws.Range("T16").Value = 0.05872177885328823

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 18.86569966666666
$ws.Range("H17").Value = 56.59709899999999
$ws.Range("I17").Value = 0.1849545940974826
$ws.Range("J17").Value = 0.1849545940974826
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.88025066666667
$ws.Range("N17").Value = 68.64075199999999
$ws.Range("O17").Value = 0.2027759882277303
$ws.Range("P17").Value = 0.2027759882277303
$ws.Range("Q17").Value = 431.651937375383
$ws.Range("R17").Value = 3884.867436378447
$ws.Range("S17").Value = 0.03750435059537577
$ws.Range("T17").Value = 0.03750435059537577

